# Normalize the "Recorded By" audit strings in column G:
#  - "System, dnasr281@gmail.com"              -> "dnasr281@gmail.com, System"
#  - "system, System, backup@backdoor.com"     -> "System, backup@backdoor.com, system"
# Every other value in the column (and the rest of the sheet) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Text

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
        $changed = $changed + 1
    }
    elseif ($val -eq "system, System, backup@backdoor.com") {
        $cell.Value = "System, backup@backdoor.com, system"
        $changed = $changed + 1
    }
}

Write-Host "Reordered Recorded-By values in $changed cell(s)."
